# Completed Excel data importer:
#  - Renamed "Business Line" -> "BusinessLine" header on BusinessLineLogins sheet
#    (this also drops the now-unused "Business Line" shared string and
#    re-points every other shared-string reference that shifted underneath it)
#  - Updated two sample amounts on the Tests sheet (Ticklers rows)
#  - Moved the remembered cell selection on the Tests / BusinessLineLogins sheets

$wb = $excel.ActiveWorkbook

# --- Tests sheet -----------------------------------------------------------
$wsTests = $wb.Worksheets.Item("Tests")

# Ticklers_DataSets sample amounts (column G, rows 10 & 11)
$wsTests.Cells.Item(10, 7).Value = 90
$wsTests.Cells.Item(11, 7).Value = 0

# --- BusinessLineLogins sheet ----------------------------------------------
$wsLogins = $wb.Worksheets.Item("BusinessLineLogins")

# Header text tweak: "Business Line" -> "BusinessLine"
$wsLogins.Range("A1").Value = "BusinessLine"

# Remembered selection moves to D32 on this (non-active) sheet
[void]$wsLogins.Range("D32").Select()

# Tests sheet stays the active tab, with its remembered selection moved to H16 -
# selected last so it remains the active/selected sheet on save.
[void]$wsTests.Range("H16").Select()
